$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D (shifts the existing D:K data right to E:L),
# adding a new first data-year column (period ending 31-Dec-2018) to each
# of the three financial statements (Income Statement, Balance Sheet,
# Cash Flow Statement).
$ws.Columns("D").Insert()

# Restore per-row number formatting in the new column D by copying the
# formats from column E (which now holds what used to be column D) across
# each of the three statement blocks.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Income Statement (new column D values) ----
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 58700
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -1000
$ws.Range("D17").Value = 12700
$ws.Range("D18").Value = 46100
$ws.Range("D20").Value = -28500
$ws.Range("D21").Value = 20300
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 17500
$ws.Range("D24").Value = 2800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 14700
$ws.Range("D27").Value = 14700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 28500
$ws.Range("D33").Value = 14700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 14700

# ---- Balance Sheet (new column D values) ----
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 24300
$ws.Range("D42").Value = 42400
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 27600
$ws.Range("D49").Value = 18300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 5200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1684800
$ws.Range("D57").Value = 1100
$ws.Range("D58").Value = 5700
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1541200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 69800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 143500
$ws.Range("D77").Value = 0

# ---- Cash Flow Statement (new column D values) ----
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 14700
$ws.Range("D83").Value = 2800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 17300
$ws.Range("D91").Value = -2000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -177700
$ws.Range("D96").Value = -8800
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 200000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 39500
